# Update the "Förändrad" (Changed) date column (C) for rows 2-140
# from serial date 45190 (2023-09-21) to 45192 (2023-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 140
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
